$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 246, shifting existing rows 246:281 down to 247:282
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row 246 with the new weekly data point.
$ws.Cells.Item(246, 1).Value = 8
$ws.Cells.Item(246, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(246, 3).Value = "Coquimbo"
$ws.Cells.Item(246, 4).Value = 45142
$ws.Cells.Item(246, 4).NumberFormat = $ws.Cells.Item(247, 4).NumberFormat
$ws.Cells.Item(246, 5).Value = 4
$ws.Cells.Item(246, 6).Value = 100112001
$ws.Cells.Item(246, 7).Value = "Berenjena"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 500
$ws.Cells.Item(246, 11).Value = 9000
$ws.Cells.Item(246, 12).Value = 10000
$ws.Cells.Item(246, 13).Value = 9500
$ws.Cells.Item(246, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(246, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(246, 16).Value = 190
$ws.Cells.Item(246, 17).Value = 50
$ws.Cells.Item(246, 18).Value = "Hortaliza"
